$wb = $excel.ActiveWorkbook

# Sheets that need the C1 cell's vertical alignment normalized to "top"
# before the value is typed in (these currently use a header style that
# lacks vertical="top", so typing the new value there first needs the
# alignment nudged so the cell re-uses the common header style).
$needsAlignFix = @("Binary Quiz 1", "Binary Quiz 2", "Binary Quiz 3")

$sheetNames = @("GroupQuiz1", "90% CI Quiz 1", "90% CI Quiz 2", "90% CI Quiz 3", "Binary Quiz 1", "Binary Quiz 2", "Binary Quiz 3")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()

    if ($needsAlignFix -contains $name) {
        $ws.Range("C1").VerticalAlignment = -4160
    }

    $ws.Range("C1").Select()
    $ws.Range("C1").Value = "Solution"
}

# The last sheet ends up with the selection one cell to the right (D1).
$lastWs = $wb.Worksheets.Item("Binary Quiz 3")
$lastWs.Range("D1").Select()

# Re-activate the first sheet, which becomes the active tab.
$wb.Worksheets.Item("GroupQuiz1").Activate()
